$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7767695188522339
$ws.Range("B1").Value = 0.9612244963645935
$ws.Range("C1").Value = 1.356578946113586
$ws.Range("D1").Value = 3.210244655609131
$ws.Range("E1").Value = 2.884607315063477
